# Com. 11/08/2025 #4: 1. Funcionamiento de ventana de resumen funcional
#
# Insert a new record row right before the existing "Puebla" row (row 9),
# pushing that row down to row 10, and populate the new row with the
# "Chiapas" record data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 9; existing row 9 (and below) shift down.
$ws.Rows.Item(9).Insert()

# Fill the new row 9 with the new record's values.
$ws.Range("A9").Value = "Rc10d1"
$ws.Range("B9").Value = "Udd529"
$ws.Range("C9").Value = "Chiapas"
$ws.Range("D9").Value = "Gastronomía"
$ws.Range("E9").Value = $true
$ws.Range("F9").Value = "adasdasd"
$ws.Range("G9").Value = "asdasdasd"
$ws.Range("H9").Value = "src/main/java/proyecto/resources/cdmx/495210901_122127526118410996_4349407357195817398_n.jpg"
$ws.Range("I9").Value = $false
